# Add a "Fig Index" column at the front of the annotation table so each row
# can be paired with its distortion-metric record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing A:J data to
# B:K and carries over column widths / cell styles automatically.
$ws.Columns("A").Insert()

# Populate the new "Fig Index" column: a header plus one synthetic key per
# row built from DOI_Mistake_FigID.
$ws.Range("A1").Value = "Fig Index"
$ws.Range("A2").Value = "10.1038:s41477-023-01482-1_log_fig1"
$ws.Range("A3").Value = "10.1038:s41477-023-01482-1_log_fig2"
$ws.Range("A4").Value = "10.1038:s41477-023-01495-w_log_fig1"
$ws.Range("A5").Value = "10.1038:s41477-023-01507-9_log_fig1"
$ws.Range("A6").Value = "10.1038:s41477-023-01539-1_log_fig1"

# Match the header formatting (bold, centered, bordered) used by the rest of
# row 1 for the whole new column, header and data alike.
$ws.Range("B1").Copy()
$ws.Range("A1:A6").PasteSpecial(-4122)

# Restore the selection to where the author left off.
$ws.Range("B2").Select() | Out-Null
